$d = $word.ActiveDocument

# Title heading and SEO title (both occurrences) -> same replacement
$d.Content.Find.Execute("Play Gypsy Moon Free: Review, Pros & Cons", $true, $false, $false, $false, $false, $true, 1, $false, "Play Gypsy Moon Slot Game for Free", 2)

# Pros list items
$d.Content.Find.Execute("Well-crafted symbols with a fortune-telling theme", $true, $false, $false, $false, $false, $true, 1, $false, "Autoplay mode and Gamble function for added excitement", 2)
$d.Content.Find.Execute("Up to 30 free spins through the Scatter symbol", $true, $false, $false, $false, $false, $true, 1, $false, "Visually appealing graphics and well-crafted symbols", 2)
$d.Content.Find.Execute("Autoplay mode and Gamble function available", $true, $false, $false, $false, $false, $true, 1, $false, "Access to up to 30 free spins for the chance to win more", 2)

# Cons list items
$d.Content.Find.Execute("Gameplay features may not be too elaborate for some players", $true, $false, $false, $false, $false, $true, 1, $false, "Gameplay may not be elaborate enough for some players", 2)
$d.Content.Find.Execute("Maximum bet per line is only €1.50", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting options with a maximum bet per line of €1.50", 2)

# SEO description
$d.Content.Find.Execute("Read our review of Gypsy Moon, a slot game with fortune-telling theme. Play for free or real money on certified online casinos. Pros and cons listed.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Gypsy Moon and play this fortune-telling themed slot game for free.", 2)
